# ESSCI template: "Update Word template to 2024"
$d = $word.ActiveDocument

# 1) Conference dates: "March 8-11, 2018" -> "March 10-13, 2024"
$d.Content.Find.Execute(" 8-11, 2018", $true, $false, $false, $false, $false, $true, 1, $false, " 10-13, 2024", 1) | Out-Null

# 2) Conference location: "Columbia, South Carolina" -> "Athens, Georgia"
$d.Content.Find.Execute("Columbia, South Carolina", $true, $false, $false, $false, $false, $true, 1, $false, "Athens, Georgia", 1) | Out-Null

# 3) Table caption run merge: ", Kelvin" + " " -> ", Kelvin "  (no visible text change, just tidy-up)
$d.Content.Find.Execute("kJoules, Kelvin ", $true, $false, $false, $false, $false, $true, 1, $false, "kJoules, Kelvin ", 1) | Out-Null

# 4) Reference citation: drop the stray "_GoBack" bookmark, merge "-" + "2626." -> "-2626."
$d.Content.Find.Execute("2619-2626.", $true, $false, $false, $false, $false, $true, 1, $false, "2619-2626.", 1) | Out-Null
